$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 25002526
$ws.Range("J28").Value = 393.42856
$ws.Range("L28").Value = 393.42856
$ws.Range("N28").Value = -1363.42856
$ws.Range("H96").Value = 970.75
$ws.Range("I96").Value = 433.2
$ws.Range("J96").Value = 1866.6666
$ws.Range("K96").Value = 1299.6
$ws.Range("L96").Value = 5599.9998
$ws.Range("M96").Value = 73.40000000000009
$ws.Range("N96").Value = -8345.9998
$ws.Range("H98").Value = 5717390
$ws.Range("I98").Value = 5883564
$ws.Range("J98").Value = 3834084
$ws.Range("K98").Value = 5883564
$ws.Range("L98").Value = 3834084
$ws.Range("M98").Value = -5882066
$ws.Range("N98").Value = -3837080
$ws.Range("H99").Value = 1127
$ws.Range("I99").Value = 720.5
$ws.Range("J99").Value = 2211
$ws.Range("K99").Value = 2161.5
$ws.Range("L99").Value = 6633
$ws.Range("M99").Value = -663.5
$ws.Range("N99").Value = -9629
$ws.Range("H101").Value = 19823.908
$ws.Range("I101").Value = 474.42856
$ws.Range("J101").Value = 28853.666
$ws.Range("K101").Value = 1423.28568
$ws.Range("L101").Value = 86560.99800000001
$ws.Range("M101").Value = 198.71432
$ws.Range("N101").Value = -89804.99800000001
$ws.Range("H104").Value = 355.42856
$ws.Range("I104").Value = 119.2
$ws.Range("K104").Value = 357.6
$ws.Range("M104").Value = 1389.4
$ws.Range("H106").Value = 47621020
$ws.Range("I106").Value = 52633550
$ws.Range("K106").Value = 52633550
$ws.Range("M106").Value = -52632919
$ws.Range("H107").Value = 911.72
$ws.Range("I107").Value = 943.1739
$ws.Range("J107").Value = 550
$ws.Range("K107").Value = 943.1739
$ws.Range("L107").Value = 550
$ws.Range("M107").Value = 976.8261
$ws.Range("N107").Value = -4390
$ws.Range("H111").Value = 7320
$ws.Range("J111").Value = 1690
$ws.Range("L111").Value = 5070
$ws.Range("N111").Value = -11204
$ws.Range("H113").Value = 9092819
$ws.Range("I113").Value = 11112889
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 11112889
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = -11109635
$ws.Range("N113").Value = -9008
$ws.Range("H118").Value = 2148
$ws.Range("J118").Value = 6066.6665
$ws.Range("L118").Value = 18199.9995
$ws.Range("N118").Value = -21513.9995
$ws.Range("H122").Value = 5717390
$ws.Range("I122").Value = 5883564
$ws.Range("J122").Value = 3834084
$ws.Range("K122").Value = 17650692
$ws.Range("L122").Value = 11502252
$ws.Range("M122").Value = -17648242
$ws.Range("N122").Value = -11507152
$ws.Range("H133").Value = 43145
$ws.Range("J133").Value = 43145
$ws.Range("L133").Value = 43145
$ws.Range("N133").Value = -53265
$ws.Range("H134").Value = 770000
$ws.Range("J134").Value = 770000
$ws.Range("L134").Value = 770000
$ws.Range("N134").Value = -780140
$ws.Range("H136").Value = 50780
$ws.Range("J136").Value = 50780
$ws.Range("L136").Value = 50780
$ws.Range("N136").Value = -60980
$ws.Range("H137").Value = 16667949
$ws.Range("I137").Value = 1328.5834
$ws.Range("J137").Value = 56667836
$ws.Range("K137").Value = 3985.7502
$ws.Range("L137").Value = 170003508
$ws.Range("M137").Value = -1435.7502
$ws.Range("N137").Value = -170008608

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 21601.941
$ws.Range("I2").Value = 22839.562
$ws.Range("J2").Value = 1800
$ws.Range("K2").Value = 22839.562
$ws.Range("L2").Value = 1800
$ws.Range("M2").Value = -22726.562
$ws.Range("N2").Value = -2026
$ws.Range("H45").Value = 626167.9
$ws.Range("I45").Value = 910303.0600000001
$ws.Range("J45").Value = 1070.4
$ws.Range("K45").Value = 910303.0600000001
$ws.Range("L45").Value = 1070.4
$ws.Range("M45").Value = -909926.0600000001
$ws.Range("N45").Value = -1824.4
$ws.Range("H116").Value = 21601.941
$ws.Range("I116").Value = 22839.562
$ws.Range("J116").Value = 1800
$ws.Range("K116").Value = 22839.562
$ws.Range("L116").Value = 1800
$ws.Range("M116").Value = -20545.562
$ws.Range("N116").Value = -6388
$ws.Range("H122").Value = 2039
$ws.Range("I122").Value = 1901.7142
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 5705.142599999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3255.142599999999
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 21601.941
$ws.Range("I3").Value = 22839.562
$ws.Range("J3").Value = 1800
$ws.Range("K3").Value = 22839.562
$ws.Range("L3").Value = 1800
$ws.Range("M3").Value = -22725.562
$ws.Range("N3").Value = -2028
$ws.Range("H99").Value = 1217
$ws.Range("I99").Value = 808.5833
$ws.Range("J99").Value = 1917.1428
$ws.Range("K99").Value = 808.5833
$ws.Range("L99").Value = 1917.1428
$ws.Range("M99").Value = 689.4167
$ws.Range("N99").Value = -4913.1428
$ws.Range("H105").Value = 1846.1538
$ws.Range("I105").Value = 1833.3334
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1833.3334
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -86.33339999999998
$ws.Range("N105").Value = -5494
$ws.Range("H107").Value = 1000541.7
$ws.Range("I107").Value = 1428941.6
$ws.Range("J107").Value = 942
$ws.Range("K107").Value = 1428941.6
$ws.Range("L107").Value = 942
$ws.Range("M107").Value = -1427021.6
$ws.Range("N107").Value = -4782

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1362159.2
$ws.Range("I31").Value = 1662.0883
$ws.Range("K31").Value = 1662.0883
$ws.Range("M31").Value = -1367.0883
$ws.Range("H34").Value = 1362159.2
$ws.Range("I34").Value = 1662.0883
$ws.Range("K34").Value = 1662.0883
$ws.Range("M34").Value = -1460.0883

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 158.25
$ws.Range("I107").Value = 109.9
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 109.9
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1810.1
$ws.Range("N107").Value = -4240
$ws.Range("H122").Value = 5971881.5
$ws.Range("I122").Value = 29342
$ws.Range("J122").Value = 16668452
$ws.Range("K122").Value = 88026
$ws.Range("L122").Value = 50005356
$ws.Range("M122").Value = -85576
$ws.Range("N122").Value = -50010256
$ws.Range("H141").Value = 40714.5
$ws.Range("J141").Value = 40714.5
$ws.Range("L141").Value = 40714.5
$ws.Range("N141").Value = -51074.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 34811.035
$ws.Range("I16").Value = 37297.035
$ws.Range("J16").Value = 1250
$ws.Range("K16").Value = 37297.035
$ws.Range("L16").Value = 1250
$ws.Range("M16").Value = -37127.035
$ws.Range("N16").Value = -1590
$ws.Range("H93").Value = 19106.334
$ws.Range("I93").Value = 4912.625
$ws.Range("J93").Value = 35327.715
$ws.Range("K93").Value = 4912.625
$ws.Range("L93").Value = 35327.715
$ws.Range("M93").Value = -3664.625
$ws.Range("N93").Value = -37823.715
$ws.Range("H122").Value = 2833.7778
$ws.Range("I122").Value = 2833.7778
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8501.3334
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6051.3334
$ws.Range("N122").ClearContents()
$ws.Range("H135").Value = 40429
$ws.Range("J135").Value = 40429
$ws.Range("L135").Value = 40429
$ws.Range("N135").Value = -50569

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 422
$ws.Range("I113").Value = 412.75
$ws.Range("J113").Value = 446.66666
$ws.Range("K113").Value = 1238.25
$ws.Range("L113").Value = 1339.99998
$ws.Range("M113").Value = 931.75
$ws.Range("N113").Value = -5679.999980000001
$ws.Range("H122").Value = 1553.619
$ws.Range("I122").Value = 1399.2222
$ws.Range("J122").Value = 2480
$ws.Range("K122").Value = 4197.6666
$ws.Range("L122").Value = 7440
$ws.Range("M122").Value = -1747.6666
$ws.Range("N122").Value = -12340
$ws.Range("H140").Value = 65733.336
$ws.Range("J140").Value = 65733.336
$ws.Range("L140").Value = 65733.336
$ws.Range("N140").Value = -76093.336
